$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country name pairs (reorder ranking) ---
$a97 = $ws.Range("A97").Value2
$a98 = $ws.Range("A98").Value2
$ws.Range("A97").Value = $a98
$ws.Range("A98").Value = $a97

$a192 = $ws.Range("A192").Value2
$a193 = $ws.Range("A193").Value2
$ws.Range("A192").Value = $a193
$ws.Range("A193").Value = $a192

$a205 = $ws.Range("A205").Value2
$a206 = $ws.Range("A206").Value2
$ws.Range("A205").Value = $a206
$ws.Range("A206").Value = $a205

# --- Update "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 11:36"

# --- Update statistic cell values ---
$ws.Range("B4").Value = 2935993
$ws.Range("C4").Value = 223
$ws.Range("D4").Value = 1260619
$ws.Range("E4").Value = 1543056
$ws.Range("B6").Value = 681251
$ws.Range("C6").Value = 6736
$ws.Range("D6").Value = 450750
$ws.Range("E6").Value = 220340
$ws.Range("G6").Value = 134
$ws.Range("H6").Value = 10161
$ws.Range("D18").Value = 181700
$ws.Range("E18").Value = 6637
$ws.Range("B46").Value = 35950
$ws.Range("C46").Value = 231
$ws.Range("D46").Value = 23746
$ws.Range("E46").Value = 10687
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 1517
$ws.Range("B60").Value = 18280
$ws.Range("C60").Value = 115
$ws.Range("D60").Value = 16615
$ws.Range("E60").Value = 959
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 706
$ws.Range("D61").Value = 10718
$ws.Range("E61").Value = 6372
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 582
$ws.Range("B65").Value = 14132
$ws.Range("C65").Value = 310
$ws.Range("D65").Value = 9410
$ws.Range("E65").Value = 4488
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 234
$ws.Range("B74").Value = 8663
$ws.Range("C74").Value = 5
$ws.Range("D74").Value = 8465
$ws.Range("E74").Value = 77
$ws.Range("B81").Value = 7253
$ws.Range("C81").Value = 5
$ws.Range("E81").Value = 224
$ws.Range("B97").Value = 4043
$ws.Range("C97").Value = 208
$ws.Range("D97").Value = 463
$ws.Range("E97").Value = 3565
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = 15
$ws.Range("B98").Value = 3969
$ws.Range("D98").Value = 914
$ws.Range("E98").Value = 3007
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 48
$ws.Range("D112").Value = 1903
$ws.Range("E112").Value = 160
$ws.Range("B119").Value = 1764
$ws.Range("C119").Value = 15
$ws.Range("E119").Value = 270
$ws.Range("B120").Value = 1700
$ws.Range("C120").Value = 21
$ws.Range("E120").Value = 205
$ws.Range("B134").Value = 1088
$ws.Range("C134").Value = 6
$ws.Range("D134").Value = 965
$ws.Range("E134").Value = 55
$ws.Range("B142").Value = 939
$ws.Range("C142").Value = 12
$ws.Range("D142").Value = 891
$ws.Range("E142").Value = 48
$ws.Range("B192").Value = 47
$ws.Range("C192").Value = 2
$ws.Range("D192").Value = 11
$ws.Range("E192").Value = 34
$ws.Range("H192").Value = 2
$ws.Range("B193").Value = 46
$ws.Range("D193").Value = 45
$ws.Range("E193").Value = 1
$ws.Range("H193").Value = 0
